$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits: "N/A" -> "NA" (and fill in the previously-blank A2) ---
$ws.Range("A2").Value = "NA"
$ws.Range("E2").Value = "NA"
$ws.Range("G2").Value = "NA"

# --- Header row formatting: build the style on a scratch cell, then copy
#     it onto the header (A1:J1) in one shot so only a single new font /
#     cellXf pair is created, matching a real "apply header style" edit. ---
$scratch = $ws.Range("Z100")
$scratch.Font.Name = "Calibri"
$scratch.HorizontalAlignment = -4108
$scratch.VerticalAlignment = -4108

$header = $ws.Range("A1:J1")
$scratch.Copy()
$header.PasteSpecial(-4122)
$scratch.Clear()

# --- Move the active selection to H2 ---
$ws.Range("H2").Select()
